# Auto-generated Excel COM-interop script
# Applies numeric corrections to columns H-N across multiple sheets
# per the Omega_Profits.xlsx diff (scheduled runner data refresh).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2008.9714
$ws.Cells.Item(19, 9).Value = 1362.1
$ws.Cells.Item(19, 10).Value = 2267.72
$ws.Cells.Item(19, 11).Value = 1362.1
$ws.Cells.Item(19, 12).Value = 2267.72
$ws.Cells.Item(19, 13).Value = -1187.1
$ws.Cells.Item(19, 14).Value = -2617.72
$ws.Cells.Item(70, 8).Value = 251059.5
$ws.Cells.Item(70, 9).Value = 1400
$ws.Cells.Item(70, 11).Value = 4200
$ws.Cells.Item(70, 13).Value = -3930
$ws.Cells.Item(73, 8).Value = 251059.5
$ws.Cells.Item(73, 9).Value = 1400
$ws.Cells.Item(73, 11).Value = 4200
$ws.Cells.Item(73, 13).Value = -3264
$ws.Cells.Item(80, 8).Value = 836
$ws.Cells.Item(80, 9).Value = 737.25
$ws.Cells.Item(80, 10).Value = 934.75
$ws.Cells.Item(80, 11).Value = 2211.75
$ws.Cells.Item(80, 12).Value = 2804.25
$ws.Cells.Item(80, 13).Value = -1213.75
$ws.Cells.Item(80, 14).Value = -4800.25
$ws.Cells.Item(83, 8).Value = 836
$ws.Cells.Item(83, 9).Value = 737.25
$ws.Cells.Item(83, 10).Value = 934.75
$ws.Cells.Item(83, 11).Value = 6635.25
$ws.Cells.Item(83, 12).Value = 8412.75
$ws.Cells.Item(83, 13).Value = -1643.25
$ws.Cells.Item(83, 14).Value = -18396.75
$ws.Cells.Item(86, 8).Value = 4316.25
$ws.Cells.Item(86, 9).Value = 4030.6
$ws.Cells.Item(86, 11).Value = 4030.6
$ws.Cells.Item(86, 13).Value = -2907.6
$ws.Cells.Item(89, 8).Value = 4316.25
$ws.Cells.Item(89, 9).Value = 4030.6
$ws.Cells.Item(89, 11).Value = 20153
$ws.Cells.Item(89, 13).Value = -14537
$ws.Cells.Item(137, 8).Value = 1857.5918
$ws.Cells.Item(137, 9).Value = 1638.325
$ws.Cells.Item(137, 11).Value = 4914.975
$ws.Cells.Item(137, 13).Value = -2364.975

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 8848.075000000001
$ws.Cells.Item(32, 9).Value = 3994.9656
$ws.Cells.Item(32, 11).Value = 3994.9656
$ws.Cells.Item(32, 13).Value = -3707.9656
$ws.Cells.Item(43, 8).Value = 23965
$ws.Cells.Item(43, 10).Value = 17000
$ws.Cells.Item(43, 12).Value = 17000
$ws.Cells.Item(43, 14).Value = -17626
$ws.Cells.Item(74, 8).Value = 2134.0908
$ws.Cells.Item(74, 9).Value = 1616
$ws.Cells.Item(74, 11).Value = 1616
$ws.Cells.Item(74, 13).Value = -742
$ws.Cells.Item(77, 8).Value = 2134.0908
$ws.Cells.Item(77, 9).Value = 1616
$ws.Cells.Item(77, 11).Value = 8080
$ws.Cells.Item(77, 13).Value = -3712

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 283.25
$ws.Cells.Item(22, 9).Value = 283.25
$ws.Cells.Item(22, 11).Value = 283.25
$ws.Cells.Item(22, 13).Value = -110.25
$ws.Cells.Item(94, 8).Value = 50000610
$ws.Cells.Item(94, 9).Value = 93750296
$ws.Cells.Item(94, 10).Value = 959.7143
$ws.Cells.Item(94, 11).Value = 93750296
$ws.Cells.Item(94, 12).Value = 959.7143
$ws.Cells.Item(94, 13).Value = -93749845
$ws.Cells.Item(94, 14).Value = -1861.7143
$ws.Cells.Item(105, 8).Value = 3659.5
$ws.Cells.Item(105, 9).Value = 2881.6667
$ws.Cells.Item(105, 10).Value = 5993
$ws.Cells.Item(105, 11).Value = 2881.6667
$ws.Cells.Item(105, 12).Value = 5993
$ws.Cells.Item(105, 13).Value = -1134.6667
$ws.Cells.Item(105, 14).Value = -9487

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 1039.4
$ws.Cells.Item(7, 9).Value = 1276.75
$ws.Cells.Item(7, 10).Value = 90
$ws.Cells.Item(7, 11).Value = 1276.75
$ws.Cells.Item(7, 12).Value = 90
$ws.Cells.Item(7, 13).Value = -1163.75
$ws.Cells.Item(7, 14).Value = -316
$ws.Cells.Item(16, 8).Value = 21227.428
$ws.Cells.Item(16, 9).Value = 273.25
$ws.Cells.Item(16, 10).Value = 49166.332
$ws.Cells.Item(16, 11).Value = 273.25
$ws.Cells.Item(16, 12).Value = 49166.332
$ws.Cells.Item(16, 13).Value = 13.75
$ws.Cells.Item(16, 14).Value = -49740.332
$ws.Cells.Item(31, 8).Value = 5316.892
$ws.Cells.Item(31, 9).Value = 4335.6665
$ws.Cells.Item(31, 11).Value = 4335.6665
$ws.Cells.Item(31, 13).Value = -4040.6665
$ws.Cells.Item(34, 8).Value = 5316.892
$ws.Cells.Item(34, 9).Value = 4335.6665
$ws.Cells.Item(34, 11).Value = 4335.6665
$ws.Cells.Item(34, 13).Value = -4133.6665
$ws.Cells.Item(94, 8).Value = 839.4545000000001
$ws.Cells.Item(94, 10).Value = 495.6
$ws.Cells.Item(94, 12).Value = 495.6
$ws.Cells.Item(94, 14).Value = -1397.6
$ws.Cells.Item(105, 8).Value = 5147.6924
$ws.Cells.Item(105, 9).Value = 3713.4443
$ws.Cells.Item(105, 11).Value = 3713.4443
$ws.Cells.Item(105, 13).Value = -1966.4443
$ws.Cells.Item(113, 8).Value = 21227.428
$ws.Cells.Item(113, 9).Value = 273.25
$ws.Cells.Item(113, 10).Value = 49166.332
$ws.Cells.Item(113, 11).Value = 273.25
$ws.Cells.Item(113, 12).Value = 49166.332
$ws.Cells.Item(113, 13).Value = 1896.75
$ws.Cells.Item(113, 14).Value = -53506.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(23, 8).Value = 322.54544
$ws.Cells.Item(23, 9).Value = 90
$ws.Cells.Item(23, 11).Value = 270
$ws.Cells.Item(23, 13).Value = -35
$ws.Cells.Item(62, 8).Value = 7140.857
$ws.Cells.Item(62, 10).Value = 7497.6665
$ws.Cells.Item(62, 12).Value = 22492.9995
$ws.Cells.Item(62, 14).Value = -23864.9995
$ws.Cells.Item(65, 8).Value = 7140.857
$ws.Cells.Item(65, 10).Value = 7497.6665
$ws.Cells.Item(65, 12).Value = 67478.9985
$ws.Cells.Item(65, 14).Value = -74342.9985
$ws.Cells.Item(92, 8).Value = 825
$ws.Cells.Item(92, 10).Value = 825
$ws.Cells.Item(92, 12).Value = 2475
$ws.Cells.Item(92, 14).Value = -4971
$ws.Cells.Item(98, 8).Value = 687.7143
$ws.Cells.Item(98, 9).Value = 450
$ws.Cells.Item(98, 11).Value = 1350
$ws.Cells.Item(98, 13).Value = 148
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 14).Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 21875500
$ws.Cells.Item(14, 10).Value = 509997.5
$ws.Cells.Item(14, 12).Value = 509997.5
$ws.Cells.Item(14, 14).Value = -510333.5
$ws.Cells.Item(18, 8).Value = 6006666.5
$ws.Cells.Item(18, 9).Value = 3000000
$ws.Cells.Item(18, 10).Value = 7509999.5
$ws.Cells.Item(18, 11).Value = 3000000
$ws.Cells.Item(18, 12).Value = 7509999.5
$ws.Cells.Item(18, 13).Value = -2999707
$ws.Cells.Item(18, 14).Value = -7510585.5
$ws.Cells.Item(19, 8).Value = 7412.5
$ws.Cells.Item(19, 10).Value = 8450
$ws.Cells.Item(19, 12).Value = 8450
$ws.Cells.Item(19, 14).Value = -9026
$ws.Cells.Item(126, 8).Value = 3459.25
$ws.Cells.Item(126, 9).Value = 1945.6666
$ws.Cells.Item(126, 11).Value = 5836.9998
$ws.Cells.Item(126, 13).Value = -3366.9998
$ws.Cells.Item(132, 8).Value = 8505.5
$ws.Cells.Item(132, 9).Value = 8505.5
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 25516.5
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 14).Value = -22986.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8584.5
$ws.Cells.Item(7, 9).Value = 8584.5
$ws.Cells.Item(7, 11).Value = 8584.5
$ws.Cells.Item(7, 13).Value = -8472.5
$ws.Cells.Item(40, 8).Value = 7899.2
$ws.Cells.Item(40, 9).Value = 8142
$ws.Cells.Item(40, 11).Value = 8142
$ws.Cells.Item(40, 13).Value = -8006
$ws.Cells.Item(122, 8).Value = 4333.3335
$ws.Cells.Item(122, 9).Value = 3000
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 9000
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -6550
$ws.Cells.Item(122, 14).Value = -19900
$ws.Cells.Item(126, 8).Value = 8584.5
$ws.Cells.Item(126, 9).Value = 8584.5
$ws.Cells.Item(126, 11).Value = 25753.5
$ws.Cells.Item(126, 13).Value = -23283.5
$ws.Cells.Item(132, 8).Value = 10691.643
$ws.Cells.Item(132, 9).Value = 35798
$ws.Cells.Item(132, 10).Value = 3844.4546
$ws.Cells.Item(132, 11).Value = 107394
$ws.Cells.Item(132, 12).Value = 11533.3638
$ws.Cells.Item(132, 13).Value = -104864
$ws.Cells.Item(132, 14).Value = -16593.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 185466.25
$ws.Cells.Item(4, 9).Value = 160514
$ws.Cells.Item(4, 11).Value = 160514
$ws.Cells.Item(4, 13).Value = -160401
$ws.Cells.Item(69, 8).Value = 24569.143
$ws.Cells.Item(69, 9).Value = 13000
$ws.Cells.Item(69, 11).Value = 13000
$ws.Cells.Item(69, 13).Value = -12251
$ws.Cells.Item(72, 8).Value = 24569.143
$ws.Cells.Item(72, 9).Value = 13000
$ws.Cells.Item(72, 11).Value = 39000
$ws.Cells.Item(72, 13).Value = -35256
$ws.Cells.Item(82, 8).Value = 74999.5
$ws.Cells.Item(82, 9).Value = 75000
$ws.Cells.Item(82, 10).Value = 74999
$ws.Cells.Item(82, 11).Value = 75000
$ws.Cells.Item(82, 12).Value = 74999
$ws.Cells.Item(82, 13).Value = -74617
$ws.Cells.Item(82, 14).Value = -75765
$ws.Cells.Item(85, 8).Value = 74999.5
$ws.Cells.Item(85, 9).Value = 75000
$ws.Cells.Item(85, 10).Value = 74999
$ws.Cells.Item(85, 11).Value = 75000
$ws.Cells.Item(85, 12).Value = 74999
$ws.Cells.Item(85, 13).Value = -73674
$ws.Cells.Item(85, 14).Value = -77651
$ws.Cells.Item(107, 8).Value = 199.4
$ws.Cells.Item(107, 9).Value = 209
$ws.Cells.Item(107, 10).Value = 173
$ws.Cells.Item(107, 11).Value = 627
$ws.Cells.Item(107, 12).Value = 519
$ws.Cells.Item(107, 13).Value = 1293
$ws.Cells.Item(107, 14).Value = -4359
$ws.Cells.Item(132, 8).Value = 11587.4
$ws.Cells.Item(132, 9).Value = 12084.25
$ws.Cells.Item(132, 11).Value = 36252.75
$ws.Cells.Item(132, 13).Value = -33722.75
